$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "A1" = 0.38057104854618728
    "A2" = -0.0099999991852968151
    "A3" = -0.0089999991700633331
    "A4" = 0.061993432799894777
    "A5" = -0.0059999991849961631
    "A6" = -0.036337352171894821
    "A7" = -0.019999999023621484
    "A8" = -0.019999999019309378
    "A9" = -0.0059999991523955742
    "A10" = -0.0059999991506956007
    "A11" = 0.050513198711037433
    "A12" = -0.0059999991475865322
    "A13" = -0.0059999991336701086
    "A14" = -0.01199999907150584
    "A15" = -0.0059999991270887065
    "A16" = -0.005999999124502331
    "A17" = -0.0059999991209593873
    "A18" = -0.0089999990911904248
    "A19" = -0.07877321503573631
    "A20" = -0.0089999991850771011
    "A21" = -0.0089999991840281623
    "A22" = -0.0089999991833966675
    "A23" = -0.071810645185122368
    "A24" = -0.041999998823588847
    "A25" = -0.041999998817220607
    "A26" = -0.005999999159961078
    "A27" = -0.0059999991578516543
    "A28" = -0.0059999991485550908
    "A29" = -0.011999999084331137
    "A30" = 0.016279072442002107
    "A31" = -0.014999999044777113
    "A32" = -0.020999998985240964
    "A33" = -0.0059999991308510303
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}
